$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.03599369093307372
$ws.Range("J2").Value = 0.03599369093307372
$ws.Range("M2").Value = 5.939783333333334
$ws.Range("N2").Value = 17.81935
$ws.Range("O2").Value = 0.371037203625045
$ws.Range("P2").Value = 0.3710372036250449
$ws.Range("Q2").Value = 0.3869432253944445
$ws.Range("R2").Value = 3.48248902855
$ws.Range("S2").Value = 0.01335499843195181
$ws.Range("T2").Value = 0.01335499843195181

# Row 3
$ws.Range("I3").Value = 0.03599369093307372
$ws.Range("J3").Value = 0.03599369093307372
$ws.Range("O3").Value = 0.06951548123833583
$ws.Range("P3").Value = 0.06951548123833583
$ws.Range("S3").Value = 0.002502118746756545
$ws.Range("T3").Value = 0.002502118746756544

# Row 4
$ws.Range("I4").Value = 0.03599369093307372
$ws.Range("J4").Value = 0.03599369093307372
$ws.Range("M4").Value = 7.898083
$ws.Range("N4").Value = 23.694249
$ws.Range("O4").Value = 0.4933652400876304
$ws.Range("P4").Value = 0.4933652400876304
$ws.Range("Q4").Value = 0.5145153516463333
$ws.Range("R4").Value = 4.630638164816999
$ws.Range("S4").Value = 0.01775803596883588
$ws.Range("T4").Value = 0.01775803596883588

# Row 5
$ws.Range("I5").Value = 0.03599369093307372
$ws.Range("J5").Value = 0.03599369093307372
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.057881
$ws.Range("N5").Value = 3.173643
$ws.Range("O5").Value = 0.06608207504898879
$ws.Range("P5").Value = 0.06608207504898879
$ws.Range("Q5").Value = 0.06891495249099999
$ws.Range("R5").Value = 0.6202345724189999
$ws.Range("S5").Value = 0.002378537785529485
$ws.Range("T5").Value = 0.002378537785529485

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.744737666666667
$ws.Range("H6").Value = 5.234213
$ws.Range("I6").Value = 0.9640063090669263
$ws.Range("J6").Value = 0.9640063090669262
$ws.Range("M6").Value = 5.939783333333334
$ws.Range("N6").Value = 17.81935
$ws.Range("O6").Value = 0.371037203625045
$ws.Range("P6").Value = 0.3710372036250449
$ws.Range("Q6").Value = 10.36336371350556
$ws.Range("R6").Value = 93.27027342155
$ws.Range("S6").Value = 0.3576822051930932
$ws.Range("T6").Value = 0.3576822051930931

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.744737666666667
$ws.Range("H7").Value = 5.234213
$ws.Range("I7").Value = 0.9640063090669263
$ws.Range("J7").Value = 0.9640063090669262
$ws.Range("O7").Value = 0.06951548123833583
$ws.Range("P7").Value = 0.06951548123833583
$ws.Range("Q7").Value = 1.941622588661667
$ws.Range("R7").Value = 17.474603297955
$ws.Range("S7").Value = 0.06701336249157928
$ws.Range("T7").Value = 0.06701336249157927

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.744737666666667
$ws.Range("H8").Value = 5.234213
$ws.Range("I8").Value = 0.9640063090669263
$ws.Range("J8").Value = 0.9640063090669262
$ws.Range("M8").Value = 7.898083
$ws.Range("N8").Value = 23.694249
$ws.Range("O8").Value = 0.4933652400876304
$ws.Range("P8").Value = 0.4933652400876304
$ws.Range("Q8").Value = 13.78008290455967
$ws.Range("R8").Value = 124.020746141037
$ws.Range("S8").Value = 0.4756072041187946
$ws.Range("T8").Value = 0.4756072041187945

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.744737666666667
$ws.Range("H9").Value = 5.234213
$ws.Range("I9").Value = 0.9640063090669263
$ws.Range("J9").Value = 0.9640063090669262
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.057881
$ws.Range("N9").Value = 3.173643
$ws.Range("O9").Value = 0.06608207504898879
$ws.Range("P9").Value = 0.06608207504898879
$ws.Range("Q9").Value = 1.845724827551
$ws.Range("R9").Value = 16.611523447959
$ws.Range("S9").Value = 0.0637035372634593
$ws.Range("T9").Value = 0.0637035372634593
